$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at row 4: Caresa / Bomba de aceite / CARESA6325 / 1
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "Caresa"
$ws.Range("B4").Value = "Bomba de aceite"
$ws.Range("C4").Value = "CARESA6325"
$ws.Range("D4").Value = 1

# 2. Swap rows 16 and 17 (Illinois Junta tapa de cilindros / Juego de tornillos)
$row16 = $ws.Range("A16:D16").Value2
$row17 = $ws.Range("A17:D17").Value2
$ws.Range("A16:D16").Value2 = $row17
$ws.Range("A17:D17").Value2 = $row16

# 3. Update quantity for JD-135-15 (row 18) from 1 to 2
$ws.Range("D18").Value = 2

# 4. Update quantity for BU-009 (row 24) from 2 to 1
$ws.Range("D24").Value = 1

# 5. Insert a new row at row 46: Nubo / Válvulas escape / 1104-EP 0.8 / 1
$ws.Rows.Item(46).Insert()
$ws.Range("A46").Value = "Nubo"
$ws.Range("B46").Value = "Válvulas escape"
$ws.Range("C46").Value = "1104-EP 0.8"
$ws.Range("D46").Value = 1
